# Aggiornamento Auto-regolazione delle colonne Excel e SQL
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (auto-adjusted) ---
# NOTE: Excel's Range.ColumnWidth (chars, Normal-style font) and the
# persisted OOXML <col width="..."> differ by a constant padding offset
# (5/6 of a character, for the default Calibri 11 workbook font), so we
# compensate to land exactly on the target stored widths.
$padding = 5 / 6
$ws.Columns("A").ColumnWidth = 11 - $padding
$ws.Columns("B").ColumnWidth = 11 - $padding
$ws.Columns("F").ColumnWidth = 32 - $padding
$ws.Columns("H").ColumnWidth = 61 - $padding

# --- Data updates ---
$data = @{
    2  = @("Elmo", "Tosi", "24/05/2002", "Millepini", "PSCNNR73L22C712C", "giancarlogolgi@example.org", "0565876678", "Incrocio Viridiana, 82 Appartamento 98`n08015, Macomer (NU)")
    3  = @("Umberto", "Veltroni", "23/06/1993", "Meri'", "LFRMRT64R20I294N", "giannuzzibenedetto@example.com", "0789106263", "Incrocio Gloria, 988 Appartamento 68`n36065, Mussolente (VI)")
    4  = @("Enrico", "Saragat", "21/02/1992", "Mottola", "MLPSTN26S12E804I", "rossettiflavia@example.com", "377242544", "Via Guariento, 77`n84014, Nocera Inferiore (SA)")
    5  = @("Piero", "Lerner", "24/06/1996", "Celzi", "PCLSVN76P44D185J", "nicola58@example.org", "+39 0824513384", "Contrada Alessia, 951 Piano 2`n13835, Botto (BI)")
    6  = @("Dolores", "Comolli", "26/09/1999", "Acquacanina", "GNNVTR03B56H634C", "graziellaroth@example.com", "0371482077", "Vicolo Federico, 28 Piano 3`n20035, Villa Cortese (MI)")
    7  = @("Armando", "Gibilisco", "11/12/1996", "San Pietro Val Lemina", "ZCHPNI99L31A766E", "guglielmosoderini@example.org", "0572513888", "Vicolo Danilo, 26 Piano 4`n74021, Carosino (TA)")
    8  = @("Galasso", "Travaglia", "06/05/2005", "Vaglie", "MRCTMT74A42B371M", "carmelo04@example.net", "0942019544", "Incrocio Bianca, 29`n28857, Santa Maria Maggiore (VB)")
    9  = @("Durante", "Brugnaro", "04/07/2004", "Valestra", "PDRSRN42H19L406I", "coriolanoovadia@example.net", "076557727", "Vicolo Ansaldo, 12`n58037, Santa Fiora (GR)")
    10 = @("Flavia", "Bonatti", "25/10/2002", "Case Di Nava", "FRNMRL92A68A373V", "lboldu@example.net", "+39 057342862", "Viale Baldassare, 98`n00079, Colle Di Fuori (RM)")
    11 = @("Raffaello", "Gigli", "17/02/2005", "San Lorenzo Pioppa", "BRRDNI61C18L453Z", "doriaclaudio@example.net", "+39 351556726", "Canale Interminelli, 7`n18021, Borgomaro (IM)")
}

# Columns "data_di_nascita" (C) and "contatto_telefonico" (G) hold values
# that Excel's automatic cell-type detection would otherwise mangle
# (leading zeros dropped from phone numbers, ambiguous dd/mm dates
# reinterpreted as mm/dd). Format those two columns as Text before typing
# the values in, then restore the default "Normal" style once the literal
# text is safely stored, so the cells end up with plain values and no
# lingering custom number format.
$ws.Range("C2:C11").NumberFormat = "@"
$ws.Range("G2:G11").NumberFormat = "@"

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($c = 1; $c -le $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}

$ws.Range("C2:C11").Style = "Normal"
$ws.Range("G2:G11").Style = "Normal"
